$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.368.41"
$ws.Range("D3").Value = "1.836.97"
$ws.Range("E3").Value = "  +3.54%  "
$ws.Range("E4").Value = "  +2.22%  "
$ws.Range("D5").Value = "'317.86"
$ws.Range("E5").Value = "  +3.52%  "
$ws.Range("D6").Value = "'1.023"
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("D7").Value = "'0.4356"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("D8").Value = "'0.3717"
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("D9").Value = "'0.07339"
$ws.Range("E9").Value = "  +2.60%  "
$ws.Range("D10").Value = "'0.8719"
$ws.Range("E10").Value = "  +3.74%  "
$ws.Range("D11").Value = "'21.34"
$ws.Range("E11").Value = "  +4.49%  "
$ws.Range("D12").Value = "1.940.99"
$ws.Range("E12").Value = "  +8.03%  "
$ws.Range("D13").Value = "'5.468"
$ws.Range("E13").Value = "  +4.14%  "
$ws.Range("D14").Value = "'6.684"
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("D15").Value = "'0.07122"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("D16").Value = "'82.12"
$ws.Range("E16").Value = "  +4.11%  "
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("D18").Value = "'0.000008976"
$ws.Range("E18").Value = "  +3.19%  "
$ws.Range("D19").Value = "'1.022"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("E20").Value = "  +3.13%  "
$ws.Range("D21").Value = "27.391.63"
$ws.Range("D22").Value = "'5.247"
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("D23").Value = "'11.15"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").Value = "2.135.61"
$ws.Range("E24").Value = "  +6.73%  "
$ws.Range("D25").Value = "'156.60"
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("D26").Value = "'1.895"
$ws.Range("E26").Value = "  +4.58%  "
$ws.Range("E27").Value = "  +2.90%  "
$ws.Range("D28").Value = "'5.238"
$ws.Range("E28").Value = "  +3.35%  "
$ws.Range("E29").Value = "  +8.60%  "
$ws.Range("D30").Value = "'115.48"
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("D31").Value = "'0.09031"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").Value = "'1.200"
$ws.Range("E32").Value = "  +7.65%  "
$ws.Range("D33").Value = "'0.7591"
$ws.Range("E33").Value = "  +4.58%  "
$ws.Range("D34").Value = "'4.459"
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("D35").Value = "'2.859"
$ws.Range("E35").Value = "  +4.12%  "
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("D37").Value = "'1.147"
$ws.Range("E37").Value = "  +4.31%  "
$ws.Range("D38").Value = "'0.01955"
$ws.Range("E38").Value = "  +3.75%  "
$ws.Range("D39").Value = "'0.05244"
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("D40").Value = "'0.5156"
$ws.Range("E40").Value = "  +4.97%  "
$ws.Range("D41").Value = "'2.776"
$ws.Range("E41").Value = "  +6.70%  "
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("D43").Value = "'6.539"
$ws.Range("E43").Value = "  +3.09%  "
$ws.Range("D44").Value = "'8.459"
$ws.Range("E44").Value = "  +6.25%  "
$ws.Range("D45").Value = "'108.24"
$ws.Range("E45").Value = "  +3.36%  "
$ws.Range("D46").Value = "'10.53"
$ws.Range("E46").Value = "  +4.12%  "
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("D48").Value = "'1.673"
$ws.Range("E48").Value = "  +2.79%  "
$ws.Range("D49").Value = "'0.4616"
$ws.Range("E49").Value = "  +3.62%  "
$ws.Range("D50").Value = "'0.06294"
$ws.Range("E50").Value = "  +1.87%  "
$ws.Range("E51").Value = "  +9.20%  "
